$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-25: coin list entries unchanged in position; update Price (D) and Volume(1h) (E).
# Values are prefixed with a leading apostrophe so Excel stores them as plain text
# (matching the workbook's original inlineStr formatting) instead of auto-converting
# numeric-looking strings (e.g. "0.9996") into real numbers.
$ws.Range("D2").Value = "'22.379.42"
$ws.Range("E2").Value = "'  +8.92%  "
$ws.Range("D3").Value = "'1.587.44"
$ws.Range("E3").Value = "'  +7.78%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "'  -0.69%  "
$ws.Range("D5").Value = "'0.9918"
$ws.Range("E5").Value = "'  +3.05%  "
$ws.Range("D6").Value = "'298.86"
$ws.Range("E6").Value = "'  +7.82%  "
$ws.Range("D7").Value = "'0.3600"
$ws.Range("E7").Value = "'  +0.28%  "
$ws.Range("D8").Value = "'0.3330"
$ws.Range("E8").Value = "'  +7.92%  "
$ws.Range("D9").Value = "'40.89"
$ws.Range("E9").Value = "'  +3.60%  "
$ws.Range("D10").Value = "'1.106"
$ws.Range("E10").Value = "'  +1.46%  "
$ws.Range("D11").Value = "'0.06893"
$ws.Range("E11").Value = "'  +4.05%  "
$ws.Range("D12").Value = "'0.9996"
$ws.Range("E12").Value = "'  -0.17%  "
$ws.Range("D13").Value = "'19.22"
$ws.Range("E13").Value = "'  +5.93%  "
$ws.Range("D14").Value = "'5.766"
$ws.Range("E14").Value = "'  +5.50%  "
$ws.Range("D15").Value = "'6.455"
$ws.Range("E15").Value = "'  +4.67%  "
$ws.Range("D16").Value = "'0.9936"
$ws.Range("E16").Value = "'  +3.28%  "
$ws.Range("D17").Value = "'0.00001058"
$ws.Range("E17").Value = "'  +3.38%  "
$ws.Range("D18").Value = "'1.589.64"
$ws.Range("E18").Value = "'  +7.98%  "
$ws.Range("D19").Value = "'0.06554"
$ws.Range("E19").Value = "'  +9.88%  "
$ws.Range("D20").Value = "'75.78"
$ws.Range("E20").Value = "'  +10.01%  "
$ws.Range("D21").Value = "'15.76"
$ws.Range("E21").Value = "'  +8.23%  "
$ws.Range("D22").Value = "'5.871"
$ws.Range("E22").Value = "'  +7.09%  "
$ws.Range("D23").Value = "'11.41"
$ws.Range("E23").Value = "'  +1.42%  "
$ws.Range("D24").Value = "'22.330.94"
$ws.Range("E24").Value = "'  +8.66%  "
$ws.Range("D25").Value = "'2.366"
$ws.Range("E25").Value = "'  +4.37%  "

# Rows 26-51: the "LEO" row was removed from the source feed, so every entry
# from the old row 27 onward shifted up by one row, and a new "Aave" entry
# was appended at the end (row 51). Update Coin (B), Link (C), Price (D) and
# Volume(1h) (E) for each of these rows to their new values.
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.470"
$ws.Range("E26").Value = "'  +17.32%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'148.56"
$ws.Range("E27").Value = "'  +3.01%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'19.01"
$ws.Range("E28").Value = "'  +10.94%  "
$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "'1.763.69"
$ws.Range("E29").Value = "'  +8.01%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'122.58"
$ws.Range("E30").Value = "'  +7.68%  "
$ws.Range("B31").Value = "HuobiToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D31").Value = "'3.926"
$ws.Range("E31").Value = "'  +0.99%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.803"
$ws.Range("E32").Value = "'  +17.68%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.9143"
$ws.Range("E33").Value = "'  +13.67%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "'0.08063"
$ws.Range("E34").Value = "'  +0.78%  "
$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.623"
$ws.Range("E35").Value = "'  +10.38%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'11.63"
$ws.Range("E36").Value = "'  +12.10%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.230"
$ws.Range("E37").Value = "'  -1.35%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.008"
$ws.Range("E38").Value = "'  +6.30%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'8.311"
$ws.Range("E39").Value = "'  +12.66%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.05954"
$ws.Range("E40").Value = "'  +2.86%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.02163"
$ws.Range("E41").Value = "'  +5.70%  "
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").Value = "'0.9916"
$ws.Range("E42").Value = "'  +3.00%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1964"
$ws.Range("E43").Value = "'  +4.61%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.5726"
$ws.Range("E44").Value = "'  +8.86%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.746"
$ws.Range("E45").Value = "'  +6.50%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.51"
$ws.Range("E46").Value = "'  +2.91%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5555"
$ws.Range("E47").Value = "'  +6.85%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'122.84"
$ws.Range("E48").Value = "'  +3.15%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.924"
$ws.Range("E49").Value = "'  +6.18%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06741"
$ws.Range("E50").Value = "'  +4.55%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'71.81"
$ws.Range("E51").Value = "'  +6.94%  "
